$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '63.780.03'
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = '  -1.19%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.413.21'
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = '  -0.06%  '

# Row 4
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = '  -0.07%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '570.27'
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = '  -0.22%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '157.90'
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = '  +0.93%  '

# Row 7
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = '  -0.02%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '3.414.11'
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = '  -0.19%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.567'
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = '  -9.57%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '7.24'
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = '  +1.28%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.119'
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = '  -3.58%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.422'
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = '  -4.48%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '3.997.15'
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = '  -0.28%  '

# Row 14
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = '  +0.26%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '27.13'
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = '  -2.41%  '

# Row 16
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = '  -8.22%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '63.838.24'
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = '  -1.12%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '3.427.98'
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = '  +0.20%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '6.09'
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = '  -3.76%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '13.63'
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = '  -2.33%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '383.98'
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = '  +1.71%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '7.79'
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = '  -2.80%  '

# Row 23
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = '  -0.08%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '71.25'
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = '  -1.61%  '

# Row 25
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = '  -6.12%  '

# Row 26
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = '  -3.78%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '9.68'
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = '  -6.63%  '

# Row 28
$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = '  +0.55%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '0.999'
$ws.Cells.Item(29, 5).NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = '  -0.24%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '6.12'
$ws.Cells.Item(30, 5).NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = '  -1.03%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.40'
$ws.Cells.Item(31, 5).NumberFormat = "@"
$ws.Cells.Item(31, 5).Value = '  -6.08%  '

# Row 32
$ws.Cells.Item(32, 5).NumberFormat = "@"
$ws.Cells.Item(32, 5).Value = '  -1.42%  '

# Row 33
$ws.Cells.Item(33, 5).NumberFormat = "@"
$ws.Cells.Item(33, 5).Value = '  +0.07%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '22.94'
$ws.Cells.Item(34, 5).NumberFormat = "@"
$ws.Cells.Item(34, 5).Value = '  -0.66%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '6.96'
$ws.Cells.Item(35, 5).NumberFormat = "@"
$ws.Cells.Item(35, 5).Value = '  -3.33%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.51'
$ws.Cells.Item(36, 5).NumberFormat = "@"
$ws.Cells.Item(36, 5).Value = '  -6.39%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '160.71'
$ws.Cells.Item(37, 5).NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = '  +0.66%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.846'
$ws.Cells.Item(38, 5).NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = '  +9.77%  '

# Row 39
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = '  -4.16%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'Hedera'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.0726'
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = '  -5.10%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'Maker'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '2.783.31'
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = '  -3.39%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '25.87'
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = '  -3.27%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '43.04'
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = '  +0.32%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '26.10'
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = '  -1.81%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '6.39'
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = '  -7.63%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '4.35'
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = '  -6.07%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.0303'
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = '  -4.44%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.40'
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = '  +9.55%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '329.35'
$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = '  +2.53%  '

# Row 50
$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = '  -4.86%  '

# Row 51
$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = '  -5.50%  '
